# sector 1 baseline regressions done
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the IV estimate / Std Error / p-value columns (B:D) for the
#    "Sector 1" regression block (rows 30-36, 39-45, 48-54) that were
#    previously left blank.
# ---------------------------------------------------------------------------
$sector1 = @{
    30 = @(-8.5996000000000006, 15.019, 0.56689999999999996)
    31 = @(-3.9148999999999998, 3.9420999999999999, 0.32069999999999999)
    32 = @(-3.3772000000000002, 15.942, 0.83220000000000005)
    33 = @(0.15570000000000001, 8.0385000000000009, 0.98450000000000004)
    34 = @(-0.55649999999999999, 0.84360000000000002, 0.50949999999999995)
    35 = @(-0.3266, 1.1733, 0.78069999999999995)
    36 = @(0.12570000000000001, 0.78559999999999997, 0.87290000000000001)
    39 = @(13.49, 14.153, 0.34050000000000002)
    40 = @(12.590999999999999, 11.974, 0.29299999999999998)
    41 = @(-5.43, 16.382999999999999, 0.74029999999999996)
    42 = @(-1.8375999999999999, 24.154, 0.93940000000000001)
    43 = @(-1.3339000000000001, 2.1318000000000001, 0.53149999999999997)
    44 = @(-0.92610000000000003, 2.0249000000000001, 0.64739999999999998)
    45 = @(-1.6332, 2.0299, 0.42109999999999997)
    48 = @(-129.26, 480.01, 0.78769999999999996)
    49 = @(-47.481000000000002, 48.655000000000001, 0.3291)
    50 = @(-24.559000000000001, 72.963999999999999, 0.73640000000000005)
    51 = @(-6.0121000000000002, 77.754000000000005, 0.93840000000000001)
    52 = @(-3.5746000000000002, 7.5008999999999997, 0.63370000000000004)
    53 = @(-1.3873, 8.4610000000000003, 0.86980000000000002)
    54 = @(-29.277000000000001, 137.35, 0.83120000000000005)
}

foreach ($row in $sector1.Keys) {
    $vals = $sector1[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}

# ---------------------------------------------------------------------------
# 2) Append a new "Sector 2" block below the existing "Sector 1" block,
#    mirroring its layout (GVC MIX / GVC BP / GVC FP sub-tables, each with
#    7 rows: transport x {F,F,F,F -> t-FE,i-FE combos} / World GVC rows),
#    but without IV estimate/Std Error/p-value numbers filled in yet.
# ---------------------------------------------------------------------------

# Sub-table headers: column A = GVC MIX / GVC FP / GVC BP (same order as the
# existing Sector 1 blocks at rows 29/38/47), column B = "Sector 2"
$ws.Cells.Item(56, 1).Value = "GVC MIX"
$ws.Cells.Item(56, 2).Value = "Sector 2"

$ws.Cells.Item(65, 1).Value = "GVC FP"
$ws.Cells.Item(65, 2).Value = "Sector 2"

$ws.Cells.Item(74, 1).Value = "GVC BP"
$ws.Cells.Item(74, 2).Value = "Sector 2"

# Each sub-table has the same 7-row layout: column A alternates "transport"
# (4x) then "World GVC" (3x); columns E/F/G are "F"/"T" flags.
$subRowsPattern = @(
    @{A = "transport"; E = "F"; F = "F"; G = "F"},
    @{A = "transport"; E = "F"; F = "F"; G = "T"},
    @{A = "transport"; E = "F"; F = "T"; G = "F"},
    @{A = "transport"; E = "T"; F = "T"; G = "F"},
    @{A = "World GVC"; E = "F"; F = "F"; G = "F"},
    @{A = "World GVC"; E = "F"; F = "F"; G = "T"},
    @{A = "World GVC"; E = "F"; F = "T"; G = "F"}
)

$blockStarts = @(57, 66, 75)
foreach ($start in $blockStarts) {
    for ($i = 0; $i -lt $subRowsPattern.Count; $i++) {
        $r = $start + $i
        $pat = $subRowsPattern[$i]
        $ws.Cells.Item($r, 1).Value = $pat.A
        $ws.Cells.Item($r, 5).Value = $pat.E
        $ws.Cells.Item($r, 6).Value = $pat.F
        $ws.Cells.Item($r, 7).Value = $pat.G
    }
}

# ---------------------------------------------------------------------------
# 3) Update the view: the user had scrolled down to / selected the newly
#    added block (mirrors the original file's selection of B30 in the
#    Sector 1 block before this edit).
# ---------------------------------------------------------------------------
$ws.Range("B57").Select()
$excel.ActiveWindow.ScrollRow = 61
